# Applies the "RBI / MIFOS strategy" update to the Repayment Schedule and
# Transactions sheets, and switches the active tab back to Repayment Schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# ---------------------------------------------------------------------------
# 1) Re-style a couple of cells in row 2 (the disbursement row).
#    F2 becomes the "amount" style (numFmt #,##0) that G2 used to have,
#    and G2 becomes a plain cell like the rest of the row.
# ---------------------------------------------------------------------------
$ws.Range("G2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("H2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F2").Value = 10000
$ws.Range("G2").ClearContents()

# E3 becomes an (empty) italic-style cell, matching the style used for the
# note cell on the Transactions sheet (font id 2 / style 12).
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("I2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null            # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("E3").ClearContents()

# ---------------------------------------------------------------------------
# 2) Clear out stray "0" placeholder values that should now be blank
#    (keeping their existing cell formatting).
# ---------------------------------------------------------------------------
$ws.Range("A2").ClearContents()
$ws.Range("D2:E2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("L2:N2").ClearContents()
$ws.Range("D4:E4").ClearContents()
$ws.Range("D5:E5").ClearContents()
$ws.Range("D6:E6").ClearContents()
$ws.Range("D7:E7").ClearContents()
$ws.Range("D8:E8").ClearContents()
$ws.Range("D9:E9").ClearContents()
$ws.Range("D10:E10").ClearContents()
$ws.Range("D11:E11").ClearContents()
$ws.Range("D12:E12").ClearContents()
$ws.Range("D13:E13").ClearContents()
$ws.Range("D14:E14").ClearContents()

# ---------------------------------------------------------------------------
# 3) Remove cells that are no longer part of the schedule (B2, O2:R2, and
#    the whole now-unused P column of running totals).
# ---------------------------------------------------------------------------
$ws.Range("B2").Clear()
$ws.Range("O2:R2").Clear()
$ws.Range("P3:P14").Clear()

# ---------------------------------------------------------------------------
# 4) Update the recalculated schedule figures (interest recalculated with a
#    slightly lower rate from period 5 onward), and move the old "P" running
#    total into "O" for rows 4-14.
# ---------------------------------------------------------------------------
$ws.Range("O4").Value = 888.49
$ws.Range("O5").Value = 888.49
$ws.Range("O6").Value = 888.49

$ws.Range("F7").Value = 820.27
$ws.Range("G7").Value = 6002.13
$ws.Range("H7").Value = 68.22
$ws.Range("O7").Value = 888.49

$ws.Range("F8").Value = 828.47
$ws.Range("G8").Value = 5173.66
$ws.Range("H8").Value = 60.02
$ws.Range("O8").Value = 888.49

$ws.Range("F9").Value = 836.75
$ws.Range("G9").Value = 4336.91
$ws.Range("H9").Value = 51.74
$ws.Range("O9").Value = 888.49

$ws.Range("F10").Value = 845.12
$ws.Range("G10").Value = 3491.79
$ws.Range("H10").Value = 43.37
$ws.Range("O10").Value = 888.49

$ws.Range("F11").Value = 853.57
$ws.Range("G11").Value = 2638.22
$ws.Range("H11").Value = 34.92
$ws.Range("O11").Value = 888.49

$ws.Range("F12").Value = 862.11
$ws.Range("G12").Value = 1776.11
$ws.Range("H12").Value = 26.38
$ws.Range("O12").Value = 888.49

$ws.Range("F13").Value = 870.73
$ws.Range("G13").Value = 905.38
$ws.Range("H13").Value = 17.76
$ws.Range("O13").Value = 888.49

$ws.Range("F14").Value = 905.38
$ws.Range("H14").Value = 9.0500000000000007
$ws.Range("K14").Value = 914.43
$ws.Range("O14").Value = 914.43

# ---------------------------------------------------------------------------
# 5) Row 2 no longer needs its custom row height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 6) Update the selection on the Repayment Schedule sheet and make it the
#    active tab again (it was Transactions before).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E20").Select() | Out-Null
